$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15419.192
$ws.Range("I32").Value = 17148.512
$ws.Range("J32").Value = 4302.143
$ws.Range("K32").Value = 17148.512
$ws.Range("L32").Value = 4302.143
$ws.Range("M32").Value = -16861.512
$ws.Range("N32").Value = -4876.143
$ws.Range("H46").Value = 14750
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 14750
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 14750
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -15388
$ws.Range("H102").Value = 1933.3334
$ws.Range("I102").Value = 1911.1111
$ws.Range("K102").Value = 1911.1111
$ws.Range("M102").Value = -289.1111000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 30000.5
$ws.Range("I22").Value = 30000.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 30000.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -29827.5
$ws.Range("N22").ClearContents()
$ws.Range("H99").Value = 1704.2693
$ws.Range("I99").Value = 1505.5555
$ws.Range("K99").Value = 1505.5555
$ws.Range("M99").Value = -7.555499999999938
$ws.Range("H105").Value = 3945.75
$ws.Range("I105").Value = 3736.9
$ws.Range("J105").Value = 4990
$ws.Range("K105").Value = 3736.9
$ws.Range("L105").Value = 4990
$ws.Range("M105").Value = -1989.9
$ws.Range("N105").Value = -8484
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 43322.332
$ws.Range("J92").Value = 43322.332
$ws.Range("L92").Value = 43322.332
$ws.Range("N92").Value = -48314.332
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 5800
$ws.Range("J64").Value = 6571.4287
$ws.Range("L64").Value = 19714.2861
$ws.Range("N64").Value = -20254.2861
$ws.Range("H67").Value = 5800
$ws.Range("J67").Value = 6571.4287
$ws.Range("L67").Value = 19714.2861
$ws.Range("N67").Value = -21586.2861
$ws.Range("H70").Value = 21338
$ws.Range("I70").Value = 30000
$ws.Range("J70").Value = 4014
$ws.Range("K70").Value = 90000
$ws.Range("L70").Value = 12042
$ws.Range("M70").Value = -89685
$ws.Range("N70").Value = -12672
$ws.Range("H73").Value = 21338
$ws.Range("I73").Value = 30000
$ws.Range("J73").Value = 4014
$ws.Range("K73").Value = 90000
$ws.Range("L73").Value = 12042
$ws.Range("M73").Value = -88908
$ws.Range("N73").Value = -14226
$ws.Range("H75").Value = 3932.6924
$ws.Range("I75").Value = 1665.2
$ws.Range("J75").Value = 5349.875
$ws.Range("K75").Value = 4995.6
$ws.Range("L75").Value = 16049.625
$ws.Range("M75").Value = -3997.6
$ws.Range("N75").Value = -18045.625
$ws.Range("H76").Value = 1650
$ws.Range("I76").Value = 1650
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 4950
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -4567
$ws.Range("N76").ClearContents()
$ws.Range("H78").Value = 3932.6924
$ws.Range("I78").Value = 1665.2
$ws.Range("J78").Value = 5349.875
$ws.Range("K78").Value = 14986.8
$ws.Range("L78").Value = 48148.875
$ws.Range("M78").Value = -9994.800000000001
$ws.Range("N78").Value = -58132.875
$ws.Range("H79").Value = 1650
$ws.Range("I79").Value = 1650
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 4950
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -3624
$ws.Range("N79").ClearContents()
$ws.Range("H92").Value = 716.8333
$ws.Range("I92").Value = 832.6667
$ws.Range("J92").Value = 601
$ws.Range("K92").Value = 2498.0001
$ws.Range("L92").Value = 1803
$ws.Range("M92").Value = -1250.0001
$ws.Range("N92").Value = -4299
$ws.Range("H94").Value = 3140
$ws.Range("I94").Value = 2870
$ws.Range("J94").Value = 3950
$ws.Range("K94").Value = 8610
$ws.Range("L94").Value = 11850
$ws.Range("M94").Value = -7934
$ws.Range("N94").Value = -13202
$ws.Range("H97").Value = 990
$ws.Range("J97").Value = 990
$ws.Range("L97").Value = 2970
$ws.Range("N97").Value = -3962
$ws.Range("H100").Value = 17007
$ws.Range("J100").Value = 6009.3335
$ws.Range("L100").Value = 18028.0005
$ws.Range("N100").Value = -19650.0005
$ws.Range("H103").Value = 1956.8572
$ws.Range("I103").Value = 424.75
$ws.Range("J103").Value = 3999.6667
$ws.Range("K103").Value = 1274.25
$ws.Range("L103").Value = 11999.0001
$ws.Range("M103").Value = -395.25
$ws.Range("N103").Value = -13757.0001
$ws.Range("H106").Value = 6666.6665
$ws.Range("J106").Value = 6666.6665
$ws.Range("L106").Value = 19999.9995
$ws.Range("N106").Value = -21891.9995
$ws.Range("H112").Value = 2739.8
$ws.Range("I112").Value = 924.75
$ws.Range("J112").Value = 10000
$ws.Range("K112").Value = 2774.25
$ws.Range("L112").Value = 30000
$ws.Range("M112").Value = -1666.25
$ws.Range("N112").Value = -32216
$ws.Range("H114").Value = 419.5
$ws.Range("I114").Value = 356
$ws.Range("J114").Value = 461.83334
$ws.Range("K114").Value = 1068
$ws.Range("L114").Value = 1385.50002
$ws.Range("M114").Value = 2186
$ws.Range("N114").Value = -7893.500019999999
$ws.Range("H121").Value = 1016.3
$ws.Range("J121").Value = 1040.3334
$ws.Range("L121").Value = 3121.0002
$ws.Range("N121").Value = -5741.0002
$ws.Range("H122").Value = 782.0263
$ws.Range("I122").Value = 490
$ws.Range("J122").Value = 1106.5
$ws.Range("K122").Value = 4410
$ws.Range("L122").Value = 9958.5
$ws.Range("M122").Value = -1960
$ws.Range("N122").Value = -14858.5
$ws.Range("H129").Value = 3334730.8
$ws.Range("I129").Value = 915
$ws.Range("J129").Value = 3847625.5
$ws.Range("K129").Value = 2745
$ws.Range("L129").Value = 11542876.5
$ws.Range("M129").Value = 2255
$ws.Range("N129").Value = -11552876.5
$ws.Range("H131").Value = 3160.1333
$ws.Range("J131").Value = 1616.0741
$ws.Range("L131").Value = 4848.2223
$ws.Range("N131").Value = -14928.2223
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 27037
$ws.Range("J92").Value = 27037
$ws.Range("L92").Value = 27037
$ws.Range("N92").Value = -30781
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1254.4546
$ws.Range("I93").Value = 1080
$ws.Range("J93").Value = 2999
$ws.Range("K93").Value = 1080
$ws.Range("L93").Value = 2999
$ws.Range("M93").Value = 168
$ws.Range("N93").Value = -5495
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 13160611
$ws.Range("I122").Value = 25002812
$ws.Range("J122").Value = 2610
$ws.Range("K122").Value = 75008436
$ws.Range("L122").Value = 7830
$ws.Range("M122").Value = -75005986
$ws.Range("N122").Value = -12730
$ws.Range("H123").Value = 41271.4
$ws.Range("J123").Value = 41271.4
$ws.Range("L123").Value = 41271.4
$ws.Range("N123").Value = -51071.4

Write-Host "Done applying Asura_Profits update."